$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update W (D6) and AD (D7) values - finalized design 5/7
$ws.Range("D6").Value = 0.18
$ws.Range("D7").Value = 9.69093

$excel.Calculate()
